$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08830399999999999
$ws.Range("H2").Value = 0.264912
$ws.Range("I2").Value = 0.04372337970871547
$ws.Range("J2").Value = 0.04372337970871546
$ws.Range("M2").Value = 0.04154133333333333
$ws.Range("N2").Value = 0.124624
$ws.Range("O2").Value = 0.002517093804502335
$ws.Range("P2").Value = 0.002517093804502335
$ws.Range("Q2").Value = 0.003668265898666666
$ws.Range("R2").Value = 0.033014393088
$ws.Range("S2").Value = 0.0001100558481767108
$ws.Range("T2").Value = 0.0001100558481767108
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08830399999999999
$ws.Range("H3").Value = 0.264912
$ws.Range("I3").Value = 0.04372337970871547
$ws.Range("J3").Value = 0.04372337970871546
$ws.Range("O3").Value = 0.8133013372545576
$ws.Range("P3").Value = 0.8133013372545578
$ws.Range("Q3").Value = 1.185257996922667
$ws.Range("R3").Value = 10.667321972304
$ws.Range("S3").Value = 0.03556028318638708
$ws.Range("T3").Value = 0.03556028318638708
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08830399999999999
$ws.Range("H4").Value = 0.264912
$ws.Range("I4").Value = 0.04372337970871547
$ws.Range("J4").Value = 0.04372337970871546
$ws.Range("O4").Value = 0.18418156894094
$ws.Range("P4").Value = 0.18418156894094
$ws.Range("Q4").Value = 0.2684154906346666
$ws.Range("R4").Value = 2.415739415712
$ws.Range("S4").Value = 0.008053040674151675
$ws.Range("T4").Value = 0.008053040674151675
$ws.Range("I5").Value = 0.5310748730197871
$ws.Range("J5").Value = 0.531074873019787
$ws.Range("M5").Value = 0.04154133333333333
$ws.Range("N5").Value = 0.124624
$ws.Range("O5").Value = 0.002517093804502335
$ws.Range("P5").Value = 0.002517093804502335
$ws.Range("Q5").Value = 0.04455565556266666
$ws.Range("R5").Value = 0.401000900064
$ws.Range("S5").Value = 0.00133676527260497
$ws.Range("T5").Value = 0.00133676527260497
$ws.Range("I6").Value = 0.5310748730197871
$ws.Range("J6").Value = 0.531074873019787
$ws.Range("O6").Value = 0.8133013372545576
$ws.Range("P6").Value = 0.8133013372545578
$ws.Range("S6").Value = 0.4319239044092872
$ws.Range("T6").Value = 0.4319239044092872
$ws.Range("I7").Value = 0.5310748730197871
$ws.Range("J7").Value = 0.531074873019787
$ws.Range("O7").Value = 0.18418156894094
$ws.Range("P7").Value = 0.18418156894094
$ws.Range("S7").Value = 0.09781420333789485
$ws.Range("T7").Value = 0.09781420333789485
$ws.Range("I8").Value = 0.4252017472714976
$ws.Range("J8").Value = 0.4252017472714976
$ws.Range("M8").Value = 0.04154133333333333
$ws.Range("N8").Value = 0.124624
$ws.Range("O8").Value = 0.002517093804502335
$ws.Range("P8").Value = 0.002517093804502335
$ws.Range("Q8").Value = 0.03567320458666667
$ws.Range("R8").Value = 0.32105884128
$ws.Range("S8").Value = 0.001070272683720654
$ws.Range("T8").Value = 0.001070272683720654
$ws.Range("I9").Value = 0.4252017472714976
$ws.Range("J9").Value = 0.4252017472714976
$ws.Range("O9").Value = 0.8133013372545576
$ws.Range("P9").Value = 0.8133013372545578
$ws.Range("S9").Value = 0.3458171496588834
$ws.Range("T9").Value = 0.3458171496588835
$ws.Range("I10").Value = 0.4252017472714976
$ws.Range("J10").Value = 0.4252017472714976
$ws.Range("O10").Value = 0.18418156894094
$ws.Range("P10").Value = 0.18418156894094
$ws.Range("S10").Value = 0.07831432492889348
$ws.Range("T10").Value = 0.07831432492889348

Write-Host "Applied 82 cell updates"
